$d = $word.ActiveDocument

# --- Part 1: highlight the "Foto de perfil" list item (yellow) ---
$rng = $d.Content.Duplicate
$rng.Find.ClearFormatting()
$found = $rng.Find.Execute("Foto de perfil", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $para = $rng.Paragraphs(1).Range
    $para.HighlightColorIndex = 7
}

# --- Part 2: "ID libri" -> "ID libro" ---
$rng2 = $d.Content.Duplicate
$rng2.Find.ClearFormatting()
$rng2.Find.Execute("ID libri", $true, $false, $false, $false, $false, $true, 1, $false, "ID libro", 2)
